$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Duplicate the existing "2021-Q4" sheet and move the copy to the end of
#    the workbook. The copy keeps the old (2021-Q4) fund data untouched and
#    becomes the new last sheet, still named "2021-Q4".
# ---------------------------------------------------------------------------
$old2021 = $wb.Worksheets.Item("2021-Q4")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$old2021.Copy($null, $lastSheet)

$newLast = $wb.Worksheets.Item($wb.Worksheets.Count)
$newLast.Name = "2021-Q4 tmp"

# ---------------------------------------------------------------------------
# 2) Rename the original sheet (still holding the old 2021-Q4 data) to
#    "2022-Q4" and overwrite its contents with the new quarter's fund table.
# ---------------------------------------------------------------------------
$old2021.Name = "2022-Q4"
$ws2 = $old2021

# Style template: the existing A2 cell already carries the "id column" style
# (s="2") used throughout this workbook for the leftmost numeric column.
$idStyleTemplate = $ws2.Range("A2")

# Make columns B-G hold plain text (fund codes / numbers-as-text), matching
# the original sheet's inline-string formatting, by forcing a text number
# format *before* the values are assigned.
$textRange = $ws2.Range("B2:G6")
$textRange.NumberFormat = "@"

# Stamp the "id column" style onto A2:A6 before writing their values.
for ($r = 2; $r -le 6; $r++) {
    $idStyleTemplate.Copy($ws2.Cells.Item($r, 1))
}

# Row 2
$ws2.Cells.Item(2,1).Value = 0
$ws2.Cells.Item(2,2).Value = "006429"
$ws2.Cells.Item(2,3).Value = "诺安恒鑫混合"
$ws2.Cells.Item(2,4).Value = "0.92"
$ws2.Cells.Item(2,5).Value = "66.36"
$ws2.Cells.Item(2,6).Value = "4.66"
$ws2.Cells.Item(2,7).Value = "0.0429"
$ws2.Cells.Item(2,8).Value = 3

# Row 3
$ws2.Cells.Item(3,1).Value = 1
$ws2.Cells.Item(3,2).Value = "008180"
$ws2.Cells.Item(3,3).Value = "同泰慧利混合A"
$ws2.Cells.Item(3,4).Value = "0.48"
$ws2.Cells.Item(3,5).Value = "93.91"
$ws2.Cells.Item(3,6).Value = "3.29"
$ws2.Cells.Item(3,7).Value = "0.0158"
$ws2.Cells.Item(3,8).Value = 10

# Row 4
$ws2.Cells.Item(4,1).Value = 2
$ws2.Cells.Item(4,2).Value = "008181"
$ws2.Cells.Item(4,3).Value = "同泰慧利混合C"
$ws2.Cells.Item(4,4).Value = "0.19"
$ws2.Cells.Item(4,5).Value = "93.91"
$ws2.Cells.Item(4,6).Value = "3.29"
$ws2.Cells.Item(4,7).Value = "0.0063"
$ws2.Cells.Item(4,8).Value = 10

# Row 5
$ws2.Cells.Item(5,1).Value = 3
$ws2.Cells.Item(5,2).Value = "001375"
$ws2.Cells.Item(5,3).Value = "金元顺安优质精选灵活配置混合C"
$ws2.Cells.Item(5,4).Value = "0.63"
$ws2.Cells.Item(5,5).Value = "68.79"
$ws2.Cells.Item(5,6).Value = "0.85"
$ws2.Cells.Item(5,7).Value = "0.0054"
$ws2.Cells.Item(5,8).Value = 4

# Row 6
$ws2.Cells.Item(6,1).Value = 4
$ws2.Cells.Item(6,2).Value = "620007"
$ws2.Cells.Item(6,3).Value = "金元顺安优质精选灵活配置混合A"
$ws2.Cells.Item(6,4).Value = "0.06"
$ws2.Cells.Item(6,5).Value = "68.79"
$ws2.Cells.Item(6,6).Value = "0.85"
$ws2.Cells.Item(6,7).Value = "0.0005"
$ws2.Cells.Item(6,8).Value = 4

# Drop the temporary text number format back to the sheet default so the
# cells end up without an explicit style index, matching the source file.
$textRange.Style = "Normal"

# ---------------------------------------------------------------------------
# 3) Finish naming the new last sheet as "2021-Q4".
# ---------------------------------------------------------------------------
$newLast.Name = "2021-Q4"

# ---------------------------------------------------------------------------
# 4) Update the "总计" (summary) sheet: existing row 2 now describes the new
#    2022-Q4 totals, and a new row 3 carries the old 2021-Q4 totals.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Push the current row 2 down into row 3 (copy formatting + values first).
$wsTotal.Cells.Item(2,1).Copy($wsTotal.Cells.Item(3,1))
$wsTotal.Cells.Item(3,1).Value = 1
$wsTotal.Cells.Item(3,2).Value = "2021-Q4"
$wsTotal.Cells.Item(3,3).Value = 2
$wsTotal.Cells.Item(3,4).Value = 0.62

# Overwrite row 2 with the new 2022-Q4 totals.
$wsTotal.Cells.Item(2,2).Value = "2022-Q4"
$wsTotal.Cells.Item(2,3).Value = 5
$wsTotal.Cells.Item(2,4).Value = 0.07

# Restore the original active sheet/tab (unchanged by the source edit).
$wsTotal.Activate()
